# Apply the "output generated at 456a3b4" gh-pages data refresh to the
# 北京-漫展信息 workbook:
#   - 展览  (sheet 1): refresh "want to go" counts (col F) + one cover image URL
#   - 演出  (sheet 2): refresh "want to go" counts (col F)
#   - 本地生活 (sheet 3): a brand-new event was scraped, so insert it as the
#     new row 2 (pushing the former row 2 down to row 3) and bump its count
#   - 全部类型 (sheet 4): same refresh as 展览, but the row numbering differs
#     slightly because this sheet also mixes in rows from other categories

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = 2620
$ws1.Range("F4").Value = 747
$ws1.Range("F6").Value = 44
$ws1.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202404/EXGNJ6RX1713250967042.jpeg"
$ws1.Range("F7").Value = 3203
$ws1.Range("F8").Value = 382
$ws1.Range("F11").Value = 8194
$ws1.Range("F12").Value = 430
$ws1.Range("F14").Value = 120
$ws1.Range("F15").Value = 56
$ws1.Range("F16").Value = 262
$ws1.Range("F17").Value = 295
$ws1.Range("F18").Value = 62
$ws1.Range("F19").Value = 326
$ws1.Range("F20").Value = 10254
$ws1.Range("F22").Value = 281
$ws1.Range("F26").Value = 166
$ws1.Range("F28").Value = 85
$ws1.Range("F29").Value = 127
$ws1.Range("F30").Value = 2647
$ws1.Range("F32").Value = 31
$ws1.Range("F35").Value = 851
$ws1.Range("F36").Value = 4044
$ws1.Range("F37").Value = 263
$ws1.Range("F38").Value = 1092
$ws1.Range("F39").Value = 2662
$ws1.Range("F42").Value = 309
$ws1.Range("F43").Value = 239
$ws1.Range("F44").Value = 29
$ws1.Range("F45").Value = 92
$ws1.Range("F46").Value = 89
$ws1.Range("F47").Value = 71
$ws1.Range("F49").Value = 57

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Performance)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F5").Value = 179
$ws2.Range("F6").Value = 40

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (Local Life) - new event scraped, insert as row 2
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Rows.Item(2).Insert()

# the row-insert doesn't fully carry over column A's bordered/bold "index"
# style onto the new row, so copy it explicitly from the row below (which
# still has the original formatting) before filling in values
$ws3.Range("A3").Copy() | Out-Null
$ws3.Range("A2").PasteSpecial(-4122) | Out-Null

$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = "'2024-04-27"
$ws3.Range("C2").Value = "北京·春日赞歌Lolita茶会门票"
$ws3.Range("D2").Value = "京密路孙河52号院 伍贰咖啡"
$ws3.Range("E2").Value = "2024.04.27 13:30-04.27 19:00"
$ws3.Range("F2").Value = 0
$ws3.Range("G2").Value = 198
$ws3.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=84317"
$ws3.Range("I2").Value = "//i2.hdslb.com/bfs/openplatform/202404/qs9qu9TO1713164732238.jpeg"

# former row 2 (塔罗集市) is now row 3; its index & "want to go" count refresh
$ws3.Range("A3").Value = 2
$ws3.Range("F3").Value = 5

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value = 2620
$ws4.Range("F3").Value = 179
$ws4.Range("F4").Value = 747
$ws4.Range("F6").Value = 44
$ws4.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202404/EXGNJ6RX1713250967042.jpeg"
$ws4.Range("F7").Value = 3203
$ws4.Range("F8").Value = 382
$ws4.Range("F11").Value = 8194
$ws4.Range("F12").Value = 430
$ws4.Range("F14").Value = 120
$ws4.Range("F15").Value = 57
$ws4.Range("F16").Value = 262
$ws4.Range("F17").Value = 295
$ws4.Range("F18").Value = 62
$ws4.Range("F19").Value = 326
$ws4.Range("F20").Value = 10255
$ws4.Range("F22").Value = 281
$ws4.Range("F26").Value = 166
$ws4.Range("F29").Value = 85
$ws4.Range("F30").Value = 127
$ws4.Range("F31").Value = 2647
$ws4.Range("F34").Value = 851
$ws4.Range("F35").Value = 4044
$ws4.Range("F36").Value = 263
$ws4.Range("F37").Value = 1095
$ws4.Range("F38").Value = 2662
$ws4.Range("F42").Value = 309
$ws4.Range("F43").Value = 239
$ws4.Range("F44").Value = 29
$ws4.Range("F45").Value = 92
$ws4.Range("F46").Value = 89
$ws4.Range("F47").Value = 71
$ws4.Range("F49").Value = 57
